# Weekly Project Update Log by Snehal
# Fill in Snehal Weekly Tasks sheet with her weekly status-log entries (cols D/E,
# plus a blank-but-present C7), size rows to fit the wrapped text, then bring this
# sheet to the front (matches the author re-activating her tab before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snehal Weekly Tasks")

$ws.Range("D2").Value = "1. Reviewed the project proposal.`n2. Had a meeting with sponsor about understanding prpject in details`n3. Participated in meeting with team and discussed about project goals"
$ws.Range("D2").WrapText = $true
$ws.Range("E2").Value = "Got some basic understanding of project"
$ws.Range("E2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 45

$ws.Range("D3").Value = "1. Brainstormed about the requirements with the team and listed the questions for sponsors`n2. Met the sponsor to clarify on those questions`n3. Worked on section 1 & 3  of Project charter`n"
$ws.Range("D3").WrapText = $true
$ws.Range("E3").Value = "1. Got clarification on the questions asked to sponsors `n2. Completed section 1 & 3 in Project Charter"
$ws.Range("E3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 75

$ws.Range("D4").Value = "1. Researched on possible libraries and software required for the project       2. Set up git account, taiga scrum board and slack communication channel and intercommunication between tools. "
$ws.Range("D4").WrapText = $true
$ws.Range("E4").Value = "1. Found that GMF is supposably to be used for this project                                                              2. All accounts and software setup is done."
$ws.Range("E4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 45

$ws.Range("D5").Value = "1. Made update on sections of project charter as per sponsor's comments.`n2. Brainstormed on Functional requirements of the project.`n3. Worked on External Interface requirements, Functional requirements, Software Quality Attributes in SRS.`n4. Tried plugin development on local to understand the technology.`n6. Researched on GMF and corresponding libraries require for project."
$ws.Range("D5").WrapText = $true
$ws.Range("E5").Value = "1. Completed changes to project charter.`n2. Completed documention External Interface Requirement, Software Quality Attribute sections in SRS.`n3. Shared undestanding with team about plugin development.`n4. Found that GMF combines both EMF and GEF and can be used for the project."
$ws.Range("E5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 120

$ws.Range("D6").Value = "1. Made changes to External Interface Requirements, Functional Requirements in SRS based on the sponsors comments.`n2. Reseached on EMF Ecore models.`n3. Started  working on design of project and worked on the section 1 & 2 of SDS."
$ws.Range("D6").WrapText = $true
$ws.Range("E6").Value = "1. Completed changes to SRS`n2. Figured out EMF Ecore model necessity.`n3. Completed the first two sections of SDS.`n4. Working on detailed design of project."
$ws.Range("E6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 75

$ws.Range("C7").Value = " "
$ws.Range("D7").Value = "1. Worked on Use Cases for each core functionalities of project                        2. Worked on Policies and Tacticis of SDS                                                              3. Shared EMF demo code with team mates and explain my understanding"
$ws.Range("D7").WrapText = $true
$ws.Range("E7").Value = "1. Use Case is integrated in SRS and some minor comments by spnosor.                                2. Completed Policies and Tactics part of SDS      3. Working on detailed design of project             4.  EMF demo is shared on git for all team members"
$ws.Range("E7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 90

# Snehal's sheet becomes the active tab/selection when the workbook was last saved
$ws.Activate() | Out-Null
$ws.Range("D8").Select() | Out-Null
